# Automatische test-sync: 2025-07-31 21:29:50
#
# Adds a new "Testmail #4" row to the Logs sheet, updates the matching
# aggregate row on the Dashboard sheet, extends the conditional formatting
# ranges on the Logs sheet, and updates the bar chart's category/value
# source ranges on the Dashboard sheet to include the new data point.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 6 with the new test mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = @(
    "Wil je 100 stuks M5-bouten bestellen?",
    "mailmind.test@zohomail.eu",
    "Testmail #4: Wil je 100 stuks M5-bouten bestellen?",
    "Bestelling / Levering",
    "Beste afzender,`nBedankt voor je e-mail. Helaas kan ik je niet helpen met het plaatsen van bestellingen via e-mail. Je kunt onze webshop bezoeken om de gewenste M5-bouten te bestellen. Mocht je nog vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam] - E-mailassistent",
    "2025-07-31 21:29:47",
    "Ja",
    "Nee",
    "Ja",
    "Nee"
)

for ($col = 1; $col -le $newRow.Length; $col++) {
    $logs.Cells.Item(6, $col).Value = $newRow[$col - 1]
}

# The multi-line "Antwoord" text causes the engine to auto-expand the row
# height; reset it back to the sheet's default so row 6 keeps the same
# (implicit, non-custom) height as the other data rows.
$logs.Rows.Item(6).AutoFit()

# ---------------------------------------------------------------------
# 2. Logs sheet: extend conditional formatting ranges from row 2:5 to 2:6
# ---------------------------------------------------------------------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "5")
    $newRange = $logs.Range($col + "2:" + $col + "6")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: append aggregate row 4 for the new category count
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(4, 2).Value = 1

# ---------------------------------------------------------------------
# 4. Dashboard sheet: extend the bar chart's category/value source ranges
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$4,'Dashboard'!`$B`$2:`$B`$4,1)"
